$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the task row "Investigar funcionamiento de Echo (framework web de Go)"
# (row 82) - the task was decided not to be pursued after all.
$ws.Rows.Item(82).Delete()

# Reflect the author's final selection/scroll position after the edit.
$ws.Range("A82").Select()
